$d = $word.ActiveDocument

# The two trailing paragraphs in this document are both empty placeholder
# paragraphs; the new content block is inserted between them, i.e. right
# after the second-to-last paragraph and right before the last one.
$count = $d.Paragraphs.Count
$anchor = $d.Paragraphs.Item($count - 1)

# Create a fresh empty paragraph right after the anchor - this is where the
# new block will live. The document's two original empty paragraphs are
# left untouched on either side of it.
$anchor.Range.InsertParagraphAfter() | Out-Null

# Re-fetch the freshly inserted (still empty) paragraph and replace its
# (empty) contents with the full OOXML for the seven new paragraphs in one
# shot, so run/proofErr/formatting markup comes through exactly as authored.
$target = $d.Paragraphs.Item($count)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:xml="http://www.w3.org/XML/1998/namespace"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>ΓxxΓxx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>zeros</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(n+</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>1,n</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>+1); </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>ΓyxΓyx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>zeros</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(n+1,1</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>);</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t> </w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>for</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> k=n+1:N </w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>phi</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> = rieurs(k:-1:k-n)</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>ΓxxΓxx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=(phi*phi'+</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ΓxxΓxx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>*(k-1-n-1</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>))/</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>(k-n-1); </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>ΓyxΓyx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t> = (phi*chants(k)+</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ΓyxΓyx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>*(k-1-n-1</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>))/</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>(k-n-1); </w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>end</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.Range.InsertXML($xml) | Out-Null
